# Shift Rota Generation - Final delivery
# Applies the updated rota assignments (Cal_Primary / Cal_Standby / BAS_Finance
# columns) and renames the BAS_FinC header to BAS_Finance.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header rename
$ws.Range("F1").Value = "BAS_Finance"

# Row 3
$ws.Range("D3").Value = "Kapil"
$ws.Range("E3").Value = "Naveen"
$ws.Range("F3").Value = "Sushvin"

# Row 4
$ws.Range("D4").Value = "Naveen"
$ws.Range("E4").Value = "Divik"
$ws.Range("F4").Value = "Kapil"

# Row 5
$ws.Range("D5").Value = "Divik"
$ws.Range("F5").Value = "Naveen"

# Row 6
$ws.Range("D6").Value = "Sushvin"

# Row 7
$ws.Range("D7").Value = "Kapil"
$ws.Range("F7").Value = "Sushvin"

# Row 8
$ws.Range("D8").Value = "Naveen"
$ws.Range("E8").Value = "Kapil"
$ws.Range("F8").Value = "Divik"

# Row 9
$ws.Range("D9").Value = "Divik"
$ws.Range("F9").Value = "Kapil"

# Row 10
$ws.Range("E10").Value = "Divik"
$ws.Range("F10").Value = "Sushvin"

# Row 11
$ws.Range("D11").Value = "Divik"
$ws.Range("E11").Value = "Sushvin"
$ws.Range("F11").Value = "Naveen"

# Row 12
$ws.Range("D12").Value = "Kapil"
$ws.Range("E12").Value = "Divik"
$ws.Range("F12").Value = "Sushvin"

# Row 13
$ws.Range("D13").Value = "Sushvin"
$ws.Range("E13").Value = "Naveen"
$ws.Range("F13").Value = "Kapil"

# Row 14
$ws.Range("D14").Value = "Divik"
$ws.Range("E14").Value = "Sushvin"
$ws.Range("F14").Value = "Naveen"

# Row 15
$ws.Range("D15").Value = "Kapil"
$ws.Range("E15").Value = "Naveen"
$ws.Range("F15").Value = "Divik"

# Row 16
$ws.Range("D16").Value = "Sushvin"
$ws.Range("E16").Value = "Divik"
$ws.Range("F16").Value = "Naveen"

# Row 17
$ws.Range("E17").Value = "Kapil"

# Row 18
$ws.Range("E18").Value = "Divik"
$ws.Range("F18").Value = "Sushvin"

# Row 19
$ws.Range("D19").Value = "Naveen"
$ws.Range("E19").Value = "Sushvin"

# Row 20
$ws.Range("D20").Value = "Divik"
$ws.Range("E20").Value = "Kapil"
$ws.Range("F20").Value = "Naveen"

# Row 21
$ws.Range("E21").Value = "Naveen"
$ws.Range("F21").Value = "Sushvin"

# Row 22
$ws.Range("F22").Value = "Kapil"
